# Apply the edit described by the diff:
#  - On the "About" sheet, add a date value (2021-04-21) in cell C1,
#    formatted with the builtin short-date number format (numFmtId 14).
#  - Minor cosmetic width tweak on the "QSfE" sheet column A (27.28515625 -> 27.33203125).

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$qsfeSheet  = $wb.Worksheets.Item("QSfE")

# Set the number format first, then the value, so the cell is stamped with
# the builtin date format (numFmtId 14) instead of a locale-inferred custom one.
$aboutSheet.Range("C1").NumberFormat = "mm-dd-yy"
$aboutSheet.Range("C1").Value = [DateTime]"2021-04-21"

# Slightly adjust column A width on the QSfE sheet to match the new value.
# (26.5 "characters" is the closest achievable ColumnWidth to the stored
# width of 27.33203125 given this engine's pixel-based quantization.)
$qsfeSheet.Columns.Item(1).ColumnWidth = 26.5
